$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous data rows (2-5) but keep the header row (row 1) intact,
# including its styling. New data (rows 2-9) will be written below.
$ws.Range("A2:T5").Clear()

# ----- Text / label columns, written column-by-column (A, then B, then C, then D) -----
# This ordering matches how the shared-strings table ends up built (new unique
# strings appended in first-seen order): FAPs(20), MuSCs(21), Adipoq(22),
# Adipor1(23), ECs(24), Resolving-Mac(25).

$colA = @("FAPs","FAPs","FAPs","FAPs","MuSCs","MuSCs","MuSCs","MuSCs")
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}

$colB = @("Adipoq","Adipoq","Adipoq","Adipoq","Adipoq","Adipoq","Adipoq","Adipoq")
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @("Adipor1","Adipor1","Adipor1","Adipor1","Adipor1","Adipor1","Adipor1","Adipor1")
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colD = @("ECs","FAPs","MuSCs","Resolving-Mac","ECs","FAPs","MuSCs","Resolving-Mac")
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

# ----- Numeric columns E..T for rows 2-9 -----
$data = @{
    2  = @(1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998, 0.603254837001998, 3, 1, 27.65020566666666, 82.95061699999999, 0.325102305838462, 0.325102305838462, 5.209446215363555, 46.885015938272, 0.1961195385175551, 0.1961195385175551)
    3  = @(1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998, 0.603254837001998, 3, 1, 14.854735, 44.564205, 0.1746572397810847, 0.1746572397810847, 2.798711299253334, 25.18840169328, 0.1053628247153571, 0.1053628247153571)
    4  = @(1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998, 0.603254837001998, 3, 1, 12.620438, 37.861314, 0.1483870877473286, 0.1483870877473286, 2.377757828202667, 21.399820453824, 0.08951522843221586, 0.08951522843221586)
    5  = @(1, 0.3333333333333333, 0.1884053333333333, 0.5652160000000001, 0.603254837001998, 0.603254837001998, 3, 1, 29.92540433333333, 89.776213, 0.3518533666331248, 0.3518533666331247, 5.638105778556445, 50.742952007008, 0.2122572453368699, 0.2122572453368699)
    6  = @(1, 0.3333333333333333, 0.1239093333333333, 0.371728, 0.396745162998002, 0.396745162998002, 3, 1, 27.65020566666666, 82.95061699999999, 0.325102305838462, 0.325102305838462, 3.426118550686222, 30.835066956176, 0.1289827673209069, 0.1289827673209069)
    7  = @(1, 0.3333333333333333, 0.1239093333333333, 0.371728, 0.396745162998002, 0.396745162998002, 3, 1, 14.854735, 44.564205, 0.1746572397810847, 0.1746572397810847, 1.840640310693333, 16.56576279624, 0.06929441506572756, 0.06929441506572755)
    8  = @(1, 0.3333333333333333, 0.1239093333333333, 0.371728, 0.396745162998002, 0.396745162998002, 3, 1, 12.620438, 37.861314, 0.1483870877473286, 0.1483870877473286, 1.563790058954667, 14.074110530592, 0.0588718593151127, 0.0588718593151127)
    9  = @(1, 0.3333333333333333, 0.1239093333333333, 0.371728, 0.396745162998002, 0.396745162998002, 3, 1, 29.92540433333333, 89.776213, 0.3518533666331248, 0.3518533666331247, 3.708036900673778, 33.372332106064, 0.1395961212962548, 0.1395961212962548)
}

for ($r = 2; $r -le 9; $r++) {
    $vals = $data[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 5).Value = $vals[$c]
    }
}
